# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 27 de Julio de 2020 a las 23:53"

# Helper: set an entire data row (A..H)
function Set-Row($r, $vals) {
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
    $ws.Cells.Item($r, 7).Value = $vals[6]
    $ws.Cells.Item($r, 8).Value = $vals[7]
}

# Simple in-place numeric updates (country stays on the same row)
Set-Row 4  @("Estados Unidos", 4423674, 51835, 2129875, 2143529, 0, 422, 150270)
Set-Row 5  @("Brasil", 2442375, 22474, 1634274, 720483, 0, 566, 87618)
Set-Row 21 @("Alemania", 207379, 638, 190600, 7574, 0, 2, 9205)
Set-Row 28 @("Egipto", 92482, 420, 34838, 52992, 0, 46, 4652)
Set-Row 33 @("Suecia", 79395, 12, 0, 0, 0, 11, 5700)
Set-Row 41 @("Israel", 63985, 2029, 27133, 36378, 0, 4, 474)

# Honduras / Barein swap position (Barein overtakes Honduras with updated numbers)
Set-Row 51 @("Barein", 39482, 351, 36110, 3231, 0, 1, 141)
Set-Row 52 @("Honduras", 39276, 838, 4922, 33238, 0, 18, 1116)

Set-Row 71 @("Costa de Marfil", 15655, 59, 10361, 5198, 0, 0, 96)
Set-Row 80 @("Estado de Palestina", 10621, 152, 3752, 6791, 0, 2, 78)
Set-Row 97 @("Republica de Yibuti", 5059, 9, 4977, 24, 0, 0, 58)

# Ruanda overtakes Islandia and Namibia with updated numbers
Set-Row 129 @("Ruanda", 1879, 58, 975, 899, 0, 0, 5)
Set-Row 130 @("Islandia", 1854, 7, 1823, 21, 0, 0, 10)
Set-Row 131 @("Namibia", 1843, 68, 101, 1734, 0, 0, 8)

Set-Row 148 @("Angola", 950, 18, 242, 667, 0, 1, 41)

# Gambia overtakes Mongolia with updated numbers
Set-Row 171 @("Gambia", 326, 49, 66, 252, 0, 2, 8)
Set-Row 172 @("Mongolia", 288, 0, 218, 70, 0, 0, 0)
